$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: selection + column width
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Columns.Item(2).ColumnWidth = 8.3
$wsSummary.Range("A5").Select() | Out-Null

# ---------------------------------------------------------------------
# Repayment schedule sheet: new Q column (accrued interest formulas)
# plus row-1 autofit + selection
# ---------------------------------------------------------------------
$wsSched = $wb.Worksheets.Item("Repayment schedule")

$wsSched.Range("Q3").Formula = "=G2*(12%/365)*B3"
$wsSched.Range("Q3").Style = "Normal"

$wsSched.Range("Q4").Formula = "=G2*(12%/365)*B4"
$wsSched.Range("Q4").Style = "Normal"

$wsSched.Range("Q5").Formula = "=G2*(12%/365)*B5"
$wsSched.Range("Q5").Style = "Normal"

$wsSched.Rows.Item(1).AutoFit() | Out-Null

$wsSched.Range("H5").Select() | Out-Null

# ---------------------------------------------------------------------
# Transactions sheet: updated transaction ids / amounts + column width
# plus selection
# ---------------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")

$wsTxn.Range("A2").Value = 1911

$wsTxn.Range("A3").Value = 1910
$wsTxn.Range("E3").Value = 100.01
$wsTxn.Range("I3").Value = 7.96

$wsTxn.Range("A4").Value = 1909
$wsTxn.Range("E4").Value = 109.78

$wsTxn.Range("A5").Value = 1905

$wsTxn.Columns.Item(1).ColumnWidth = 4.2

$wsTxn.Range("I4").Select() | Out-Null

# ---------------------------------------------------------------------
# NewLoanInput becomes the active tab (selected last so it "wins")
# ---------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("NewLoanInput")
$wsInput.Activate() | Out-Null
